$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.346.35'
$ws.Range('E2').Value = '  -3.36%  '
$ws.Range('D3').Value = '2.935.44'
$ws.Range('E3').Value = '  -4.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '493.38'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -6.92%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '132.85'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -7.39%  '
$ws.Range('E7').Value = '  +0.09%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.422'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -5.92%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '7.09'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -6.79%  '
$ws.Range('E10').Value = '  -7.47%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.349'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -5.94%  '
$ws.Range('D12').Value = '3.443.66'
$ws.Range('E12').Value = '  -4.30%  '
$ws.Range('E13').Value = '  -3.83%  '
$ws.Range('E14').Value = '  -5.73%  '
$ws.Range('E15').Value = '  -10.17%  '
$ws.Range('D16').Value = '56.452.26'
$ws.Range('E16').Value = '  -3.19%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '5.94'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('D18').Value = '2.937.07'
$ws.Range('E18').Value = '  -4.51%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.40'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -5.86%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.69'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -6.31%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '315.33'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -7.97%  '
$ws.Range('E22').Value = '  -0.17%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.73'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.37%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.479'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -5.35%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '62.36'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -4.74%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -4.96%  '
$ws.Range('D28').Value = '0.0₃0847'
$ws.Range('E28').Value = '  -13.20%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.39'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -8.90%  '
$ws.Range('E30').Value = '  -7.30%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.74'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -7.24%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '19.84'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -6.53%  '
$ws.Range('E33').Value = '  -9.71%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '150.84'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('E35').Value = '  -8.66%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.63'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -6.06%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.19'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -10.02%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '23.56'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -10.65%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0647'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -7.82%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '37.27'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').Value = '2.966.35'
$ws.Range('E41').Value = '  -4.56%  '
$ws.Range('E42').Value = '  +0.00%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.65'
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.632'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -5.58%  '
$ws.Range('D45').Value = '2.127.75'
$ws.Range('E45').Value = '  -9.13%  '
$ws.Range('E46').Value = '  -10.09%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '5.81'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -4.31%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.903'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -13.25%  '
$ws.Range('E49').Value = '  -6.56%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '18.76'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -7.25%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0842'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -6.95%  '
